# Rewrite excel bulk spec parser test fixture:
#  - sheet "Should pass - human readable row will be duplicated" is renamed
#    to "Should pass - human readable row will be ignored"
#  - the old "3 cols, 4 data" sheet is dropped; instead the
#    "3 cols, 2 spec IDs, header dup error" sheet (whose first four rows
#    already match the old "3 cols, 4 data" header/data rows) is renamed to
#    "3 cols, 4 data" and extended with the extra header cell + data rows,
#    becoming the active sheet with H7 selected.

$wb = $excel.ActiveWorkbook

# --- rename the first sheet -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
try {
    $ws1.Name = "Should pass - human readable row will be ignored"
} catch {
    Write-Output "rename failed: $_"
}

# --- drop the old "3 cols, 4 data" sheet; its content is recreated below ---
$oldDataSheet = $wb.Worksheets.Item("3 cols, 4 data")
$oldDataSheet.Delete()

# --- repurpose the "header dup error" sheet as the new "3 cols, 4 data" ---
$ws = $wb.Worksheets.Item("3 cols, 2 spec IDs, header dup error")
$ws.Name = "3 cols, 4 data"

# Row 2 gains a third header cell (reuses the shared "head3" string already
# used on the other sheets).
$ws.Range("C2").Value = "head3"

# Rows 5 and 6: plain numeric data rows (same as the old "3 cols, 4 data"
# sheet used to contain).
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1

# New row 7: four numeric data cells plus an explanatory note in H7 so the
# used range grows to A1:H7 (columns 5-8 of row 5 stay blank -> padded NaNs).
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("H7").Value = "This will make the dimensions 8x6, which will cause row 5 to be padded with NaNs in positions 5-8. Tests that the right row length is reported and the NaNs are ignored."

# This sheet remains the active tab with H7 as the selected cell.
$ws.Activate()
$ws.Range("H7").Select()
